$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/15/2025  Through  9/21/2025"

# --- Style-flip cells: use Copy() from a donor cell with the right target style+type ---
# Donor for "N/A" text cells (style 13, shared string "0"): D15
# Donor for numeric style 14: D23 (already numeric, style 14)
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("D27").Copy($ws.Range("C27"))
$ws.Range("D27").Copy($ws.Range("C28"))
$ws.Range("D27").Copy($ws.Range("F31"))
$ws.Range("D27").Copy($ws.Range("C33"))

# C23 goes from "N/A" text to numeric 1 with style 14 - copy format from D23 then set value
$ws.Range("D23").Copy()
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value = 1

# --- Bulk numeric value updates ---
# Row 15
$ws.Range("N15").Value = -7.142857142857

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -6.25
$ws.Range("I16").Value = 168
$ws.Range("J16").Value = 178
$ws.Range("K16").Value = -5.617977528089
$ws.Range("L16").Value = 41.176470588235
$ws.Range("M16").Value = 31.25
$ws.Range("N16").Value = -81.974248927038

# Row 17
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 39
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 21.875
$ws.Range("I17").Value = 281
$ws.Range("J17").Value = 254
$ws.Range("K17").Value = 10.629921259842
$ws.Range("L17").Value = 44.102564102564
$ws.Range("M17").Value = 264.935064935065
$ws.Range("N17").Value = -9.061488673139

# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 162
$ws.Range("J18").Value = 160
$ws.Range("K18").Value = 1.25
$ws.Range("L18").Value = 5.194805194805
$ws.Range("M18").Value = 95.180722891566
$ws.Range("N18").Value = -72.954924874791

# Row 19
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -29.411764705882
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 66
$ws.Range("H19").Value = -13.636363636363
$ws.Range("I19").Value = 495
$ws.Range("J19").Value = 567
$ws.Range("K19").Value = -12.698412698412
$ws.Range("L19").Value = -3.883495145631
$ws.Range("M19").Value = 39.044943820224
$ws.Range("N19").Value = -37.814070351758

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 36.363636363636
$ws.Range("I20").Value = 84
$ws.Range("J20").Value = 65
$ws.Range("K20").Value = 29.230769230769
$ws.Range("L20").Value = 50
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -84.210526315789

# Row 21
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -2.702702702702
$ws.Range("F21").Value = 145
$ws.Range("G21").Value = 144
$ws.Range("H21").Value = 0.694444444444
$ws.Range("I21").Value = 1205
$ws.Range("J21").Value = 1233
$ws.Range("K21").Value = -2.270884022708
$ws.Range("L21").Value = 15.642994241842
$ws.Range("M21").Value = 74.891146589259
$ws.Range("N21").Value = -62.249373433584

# Row 22
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -16.666666666666
$ws.Range("I22").Value = 45
$ws.Range("J22").Value = 38
$ws.Range("K22").Value = 18.421052631578
$ws.Range("L22").Value = -4.255319148936
$ws.Range("M22").Value = -6.25

# Row 23
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -62.5
$ws.Range("I23").Value = 27
$ws.Range("J23").Value = 38
$ws.Range("K23").Value = -28.947368421052
$ws.Range("L23").Value = -28.947368421052
$ws.Range("M23").Value = 50

# Row 24
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = 5
$ws.Range("F24").Value = 147
$ws.Range("G24").Value = 151
$ws.Range("H24").Value = -2.649006622516
$ws.Range("I24").Value = 1243
$ws.Range("J24").Value = 1630
$ws.Range("K24").Value = -23.742331288343
$ws.Range("L24").Value = -19.075520833333
$ws.Range("M24").Value = 12.184115523465

# Row 25
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = 8.333333333333
$ws.Range("F25").Value = 96
$ws.Range("G25").Value = 119
$ws.Range("H25").Value = -19.327731092437
$ws.Range("I25").Value = 960
$ws.Range("J25").Value = 1472
$ws.Range("K25").Value = -34.782608695652
$ws.Range("L25").Value = -30.434782608695

# Row 26
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 100
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 16.666666666666
$ws.Range("I26").Value = 317
$ws.Range("J26").Value = 332
$ws.Range("K26").Value = -4.518072289156
$ws.Range("L26").Value = 8.934707903780
$ws.Range("M26").Value = 25.296442687747

# Row 27
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 0

# Row 28
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 50
$ws.Range("J28").Value = 43
$ws.Range("K28").Value = 34.883720930232
$ws.Range("L28").Value = 28.888888888888

# Row 31

# Row 33
